$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder columns: move column C (State) to become column A.
# Address shifts from A->B, City shifts from B->C. Zip (D) and BPO (E) stay put.
$ws.Columns.Item(3).Cut()
$ws.Columns.Item(1).Insert()

# The cut left column C (now vacated/duplicated) with stale width metadata;
# clear its formatting so no stray column-width definition remains.
$ws.Columns.Item(3).ClearFormats()

# Update the active cell/selection on the sheet (cosmetic, matches the saved view state).
$ws.Range("H15").Select()

Write-Host "Column reordering complete"
